$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-13 Tuesday" "2026-01-14 Wednesday"

Replace-Text "30÷3=" "44÷5="
Replace-Text "44÷4=" "69÷8="
Replace-Text "60÷7=" "52÷7="
Replace-Text "79÷5=" "96÷7="
Replace-Text "56÷5=" "21÷4="

Replace-Text "98÷5=" "12÷6="
Replace-Text "84÷3=" "49÷8="
Replace-Text "97÷2=" "45÷7="
Replace-Text "70÷2=" "96÷4="
Replace-Text "80÷8=" "41÷4="

Replace-Text "81÷2=" "35÷8="
Replace-Text "48÷9=" "86÷8="
Replace-Text "33÷8=" "19÷5="
Replace-Text "66÷2=" "54÷5="
Replace-Text "35÷3=" "78÷5="

Replace-Text "57÷2=" "87÷7="
Replace-Text "20÷7=" "98÷8="
Replace-Text "51÷6=" "81÷5="
Replace-Text "40÷9=" "23÷8="
Replace-Text "93÷5=" "88÷9="

Replace-Text "80÷2=" "51÷2="
Replace-Text "34÷5=" "34÷9="
Replace-Text "24÷9=" "83÷4="
Replace-Text "36÷5=" "63÷8="
Replace-Text "96÷6=" "33÷5="
